$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("B2").Value = 50000
$ws.Range("D2").Value = 0.4806733659468704
$ws.Range("E2").Value = 3.728306393057234
$ws.Range("F2").Value = 0.9524
$ws.Range("H2").Value = 3.731062336476929

# Update row 3 values
$ws.Range("B3").Value = 50000
$ws.Range("D3").Value = 0.7652061428218204
$ws.Range("E3").Value = 3.605889614447373
$ws.Range("F3").Value = 1.4484
$ws.Range("H3").Value = 3.731062336476929

# Delete rows 4 and 5 (previously rows for #PCs 3 and 4)
$ws.Range("A4:H5").Delete()
